$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update agenda content (Subject / Content columns) ---
$ws.Range("B4").Value = "Basic Statistics"

$ws.Range("B5").Value = "Regression Modeling 1"

$ws.Range("B6").Value = "Regression Modeling 2"
$ws.Range("C6").Value = "e.g., stepwise regression, visualisation, tables"

$ws.Range("B7").Value = "Regression Modeling 3"
$ws.Range("C7").Value = "Continues (may include multilevel and some ANOVAs)`$^a`$"

$ws.Range("B8").Value = "Factor Analysis"
$ws.Range("C8").Value = "Confirmatory Factor Analysis and Structural Equation Modeling`$^a`$"

# --- Widen column B to fit the longer subject/content text ---
# (62.0 is the input value that the engine's character->pixel->character
#  quantization snaps closest to the authored stored width of 62.85546875)
$ws.Range("B1").ColumnWidth = 62

# --- Move the active selection to reflect where the author left off editing ---
$ws.Range("B10").Select()
